$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1.76
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 1.95
$ws.Range("Q2").Value = 1.82
$ws.Range("T2").Value = 1.76
$ws.Range("W2").Value = 1.24
$ws.Range("AK2").Value = 70
$ws.Range("AL2").Value = 140
$ws.Range("AM2").Value = 200
$ws.Range("AN2").Value = 75
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 9.4
$ws.Range("H3").Value = 1.41
$ws.Range("I3").Value = 1.52
$ws.Range("J3").Value = 4.9
$ws.Range("K3").Value = 5.9
$ws.Range("L3").Value = 1.29
$ws.Range("N3").Value = 4.5
$ws.Range("O3").Value = 1.21
$ws.Range("P3").Value = 2.26
$ws.Range("Q3").Value = 1.6
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 2.46
$ws.Range("T3").Value = 1.84
$ws.Range("U3").Value = 1.92
$ws.Range("V3").Value = 2.92
$ws.Range("W3").Value = 1.12
$ws.Range("X3").Value = 46
$ws.Range("Y3").Value = 22
$ws.Range("Z3").Value = 21
$ws.Range("AC3").Value = 42
$ws.Range("AD3").Value = 40
$ws.Range("AE3").Value = 70
$ws.Range("AO3").Value = 15
$ws.Range("M4").Value = 1.07
$ws.Range("P4").Value = 1.85
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 4.2
$ws.Range("Q5").Value = 2.04
$ws.Range("AC5").Value = 42
$ws.Range("AN5").Value = 85
$ws.Range("P6").Value = 1.8
$ws.Range("T6").Value = 1.73
$ws.Range("U6").Value = 2.02
$ws.Range("V6").Value = 1.37
$ws.Range("AJ6").Value = 85
$ws.Range("J7").Value = 3.9
$ws.Range("K7").Value = 5.3
$ws.Range("O7").Value = 1.3
$ws.Range("F8").Value = 2.08
$ws.Range("G8").Value = 2.32
$ws.Range("H8").Value = 3.7
$ws.Range("I8").Value = 4.9
$ws.Range("J8").Value = 3.1
$ws.Range("K8").Value = 3.95
$ws.Range("L8").Value = 1.5
$ws.Range("N8").Value = 2.92
$ws.Range("O8").Value = 1.38
$ws.Range("P8").Value = 1.67
$ws.Range("Q8").Value = 2.18
$ws.Range("S8").Value = 4
$ws.Range("T8").Value = 1.87
$ws.Range("U8").Value = 1.86
$ws.Range("V8").Value = 1.25
$ws.Range("W8").Value = 1.75
$ws.Range("Y8").Value = 14
$ws.Range("Z8").Value = 32
$ws.Range("AB8").Value = 8.800000000000001
$ws.Range("AC8").Value = 8.6
$ws.Range("AD8").Value = 19
$ws.Range("AF8").Value = 14
$ws.Range("AH8").Value = 22
$ws.Range("AJ8").Value = 30
$ws.Range("AK8").Value = 29
$ws.Range("AN8").Value = 24
$ws.Range("AO8").Value = 600
$ws.Range("F9").Value = 2.48
$ws.Range("I9").Value = 3.1
$ws.Range("J9").Value = 3.3
$ws.Range("K9").Value = 3.8
$ws.Range("N9").Value = 3.6
$ws.Range("P9").Value = 1.89
$ws.Range("Q9").Value = 1.83
$ws.Range("T9").Value = 1.67
$ws.Range("U9").Value = 2.1
$ws.Range("V9").Value = 1.47
$ws.Range("X9").Value = 17
$ws.Range("Y9").Value = 13.5
$ws.Range("Z9").Value = 21
$ws.Range("AA9").Value = 280
$ws.Range("AB9").Value = 12.5
$ws.Range("AC9").Value = 8.800000000000001
$ws.Range("AD9").Value = 13.5
$ws.Range("AE9").Value = 85
$ws.Range("AF9").Value = 36
$ws.Range("AG9").Value = 13.5
$ws.Range("AH9").Value = 29
$ws.Range("AJ9").Value = 170
$ws.Range("AK9").Value = 32
$ws.Range("AL9").Value = 55
$ws.Range("AN9").Value = 44
$ws.Range("AO9").Value = 46
$ws.Range("P11").Value = 2.88
$ws.Range("AG11").Value = 10.5
$ws.Range("AL11").Value = 25
$ws.Range("AM11").Value = 580
$ws.Range("AN11").Value = 5.5
$ws.Range("G12").Value = 2.88
$ws.Range("I12").Value = 2.8
$ws.Range("J12").Value = 3.5
$ws.Range("K12").Value = 3.8
$ws.Range("P12").Value = 2.24
$ws.Range("R12").Value = 1.5
$ws.Range("S12").Value = 2.68
$ws.Range("U12").Value = 2.46
$ws.Range("W12").Value = 1.53
$ws.Range("AK12").Value = 75
$ws.Range("AL12").Value = 95
$ws.Range("I13").Value = 2.18
$ws.Range("K13").Value = 4
$ws.Range("R13").Value = 1.54
$ws.Range("T13").Value = 1.6
$ws.Range("AF13").Value = 28
$ws.Range("AH13").Value = 15.5
$ws.Range("AJ13").Value = 160
$ws.Range("F14").Value = 6.4
$ws.Range("H14").Value = 1.59
$ws.Range("O14").Value = 1.24
$ws.Range("P14").Value = 2.16
$ws.Range("S14").Value = 2.92
$ws.Range("U14").Value = 2.06
$ws.Range("X14").Value = 17.5
$ws.Range("AC14").Value = 9.6
$ws.Range("AE14").Value = 16
$ws.Range("AJ14").Value = 200
$ws.Range("AM14").Value = 130
$ws.Range("AO14").Value = 8
$ws.Range("K15").Value = 950
